# Update the LR-pairs worksheet with newly-recomputed TPM values.
# The sending-cluster x target-cluster combinations are now fully populated
# (4 sending clusters x 3 target clusters = 12 data rows instead of 8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2..13, columns A..T (20 columns).
$data = @(
    @("ECs","Gnai2","Agtr2","ECs",3,1,169.915657,509.746971,0.4441184931734509,0.4441184931734509,1,0.3333333333333333,0.045282,0.135846,0.001264082841858775,0.001264082841858776,7.694120780274,69.247087022466,0.0005614025669727329,0.000561402566972733),
    @("ECs","Gnai2","Agtr2","FAPs",3,1,169.915657,509.746971,0.4441184931734509,0.4441184931734509,3,1,34.682839,104.048517,0.9681988800594137,0.9681988800594138,5893.157375310224,53038.41637779201,0.4299950277042095,0.4299950277042096),
    @("ECs","Gnai2","Agtr2","MuSCs",3,1,169.915657,509.746971,0.4441184931734509,0.4441184931734509,2,0.6666666666666666,1.093898333333333,3.281695,0.03053703709872749,0.03053703709872749,185.8704539995384,1672.834085995845,0.01356206290226862,0.01356206290226862),
    @("FAPs","Gnai2","Agtr2","ECs",3,1,68.382243,205.146729,0.1787346690539575,0.1787346690539575,1,0.3333333333333333,0.045282,0.135846,0.001264082841858775,0.001264082841858776,3.096484727526,27.868362547734,0.0002259354283964143,0.0002259354283964144),
    @("FAPs","Gnai2","Agtr2","FAPs",3,1,68.382243,205.146729,0.1787346690539575,0.1787346690539575,3,1,34.682839,104.048517,0.9681988800594137,0.9681988800594138,2371.690324427877,21345.21291985089,0.1730507064058316,0.1730507064058316),
    @("FAPs","Gnai2","Agtr2","MuSCs",3,1,68.382243,205.146729,0.1787346690539575,0.1787346690539575,2,0.6666666666666666,1.093898333333333,3.281695,0.03053703709872749,0.03053703709872749,74.803221647295,673.228994825655,0.00545802721972948,0.00545802721972948),
    @("MuSCs","Gnai2","Agtr2","ECs",3,1,53.27463399999999,159.823902,0.1392470275793777,0.1392470275793778,1,0.3333333333333333,0.045282,0.135846,0.001264082841858775,0.001264082841858776,2.412381976787999,21.711437791092,0.0001760197783429271,0.0001760197783429271),
    @("MuSCs","Gnai2","Agtr2","FAPs",3,1,53.27463399999999,159.823902,0.1392470275793777,0.1392470275793778,3,1,34.682839,104.048517,0.9681988800594137,0.9681988800594138,1847.715553805926,16629.43998425333,0.1348188161539558,0.1348188161539559),
    @("MuSCs","Gnai2","Agtr2","MuSCs",3,1,53.27463399999999,159.823902,0.1392470275793777,0.1392470275793778,2,0.6666666666666666,1.093898333333333,3.281695,0.03053703709872749,0.03053703709872749,58.27703334154333,524.4933000738899,0.004252191647078987,0.004252191647078988),
    @("Resolving-Mac","Gnai2","Agtr2","ECs",3,1,91.01828266666666,273.054848,0.2378998101932138,0.2378998101932138,1,0.3333333333333333,0.045282,0.135846,0.001264082841858775,0.001264082841858776,4.121489875711999,37.093408881408,0.000300725068146701,0.000300725068146701),
    @("Resolving-Mac","Gnai2","Agtr2","FAPs",3,1,91.01828266666666,273.054848,0.2378998101932138,0.2378998101932138,3,1,34.682839,104.048517,0.9681988800594137,0.9681988800594138,3156.772443784491,28410.95199406042,0.2303343297954167,0.2303343297954167),
    @("Resolving-Mac","Gnai2","Agtr2","MuSCs",3,1,91.01828266666666,273.054848,0.2378998101932138,0.2378998101932138,2,0.6666666666666666,1.093898333333333,3.281695,0.03053703709872749,0.03053703709872749,99.56474771192889,896.08272940736,0.007264755329650398,0.007264755329650398)
)

$startRow = 2
$numNewRows = $data.Count              # 12
$numOldRows = 8                        # previously rows 2..9 held data
$lastCol = 20                          # column T

# Remove any leftover rows beyond the new data block (old sheet had rows 2..9;
# new sheet has rows 2..13, so nothing trails off the end here, but clear
# defensively in case old row count was ever larger than the new one).
if ($numOldRows -gt $numNewRows) {
    $extraStart = $startRow + $numNewRows
    $extraEnd = $startRow + $numOldRows - 1
    $ws.Range($ws.Cells.Item($extraStart, 1), $ws.Cells.Item($extraEnd, $lastCol)).ClearContents() | Out-Null
}

# Write the data rows
for ($r = 0; $r -lt $data.Count; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($startRow + $r, $c + 1).Value2 = $row[$c]
    }
}
